$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same date value (45172 -> 2023-09-03)
# for every data row (2..439). Update it to 45175 (2023-09-06).
$ws.Range("C2:C439").Value = 45175
